$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.384.62"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.876.68"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7121"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3124"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08454"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.47%  "

$ws.Range("D12").Value = "1.871.39"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "29.386.14"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008248"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "2.119.76"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.793"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1597"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.073"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.00%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.430"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.334"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.290"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05292"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.942"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7445"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -14.38%  "

$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("D39").Value = "1.225.62"
$ws.Range("E39").Value = "  +4.61%  "

$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.484"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8943"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "2.016.61"
$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("E47").Value = "  +1.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5212"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("E49").Value = "  +4.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.402"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4331"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
